$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 value
$ws.Range("B2").Value = 66

# Update row 3: label and value
$ws.Range("A3").Value = "Wrong_Entity_Event_as_NonEvent"
$ws.Range("B3").Value = 48

# Update row 4: label and value
$ws.Range("A4").Value = "Correct"
$ws.Range("B4").Value = 37

# Update row 5 value
$ws.Range("B5").Value = 2

# Update row 6: label only (value unchanged)
$ws.Range("A6").Value = "False_I-NonEvent"
